$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.048.64"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "1.800.55"
$ws.Range("E3").Value = "  -1.94%  "

$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +1.01%  "

$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").Value = "307.07"
$ws.Range("E6").Value = "  -2.31%  "

$ws.Range("D7").Value = "0.4240"
$ws.Range("E7").Value = "  -1.17%  "

$ws.Range("D8").Value = "0.3617"
$ws.Range("E8").Value = "  -1.25%  "

$ws.Range("D9").Value = "0.07214"
$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").Value = "0.8437"
$ws.Range("E10").Value = "  -3.24%  "

$ws.Range("D11").Value = "20.27"
$ws.Range("E11").Value = "  -2.20%  "

$ws.Range("D12").Value = "1.889.50"
$ws.Range("E12").Value = "  +5.67%  "

$ws.Range("E13").Value = "  -3.05%  "

$ws.Range("D14").Value = "6.391"
$ws.Range("E14").Value = "  -2.34%  "

$ws.Range("D15").Value = "0.06823"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("D17").Value = "80.50"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").Value = "0.000008709"
$ws.Range("E18").Value = "  -2.51%  "

$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("D20").Value = "14.98"
$ws.Range("E20").Value = "  -2.93%  "

$ws.Range("D21").Value = "27.302.71"
$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("D22").Value = "5.064"
$ws.Range("E22").Value = "  -2.10%  "

$ws.Range("D23").Value = "11.08"
$ws.Range("E23").Value = "  +1.82%  "

$ws.Range("D24").Value = "2.079.61"
$ws.Range("E24").Value = "  +3.42%  "

$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("D26").Value = "153.37"
$ws.Range("E26").Value = "  -0.88%  "

$ws.Range("D27").Value = "18.29"
$ws.Range("E27").Value = "  -3.39%  "

$ws.Range("D28").Value = "115.14"
$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").Value = "5.037"
$ws.Range("E29").Value = "  -2.83%  "

$ws.Range("E30").Value = "  -12.00%  "

$ws.Range("D31").Value = "0.08944"
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("D32").Value = "0.7232"
$ws.Range("E32").Value = "  -4.84%  "

$ws.Range("D33").Value = "2.853"
$ws.Range("E33").Value = "  -3.74%  "

$ws.Range("D34").Value = "4.335"
$ws.Range("E34").Value = "  -4.72%  "

$ws.Range("D35").Value = "1.090"
$ws.Range("E35").Value = "  -4.74%  "

$ws.Range("E36").Value = "  +0.92%  "

$ws.Range("D37").Value = "1.085"
$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("D38").Value = "0.01897"
$ws.Range("E38").Value = "  -2.59%  "

$ws.Range("D39").Value = "0.05090"
$ws.Range("E39").Value = "  -4.41%  "

$ws.Range("D40").Value = "0.4972"
$ws.Range("E40").Value = "  -2.56%  "

$ws.Range("D41").Value = "0.1622"
$ws.Range("E41").Value = "  -3.07%  "

$ws.Range("D42").Value = "2.535"
$ws.Range("E42").Value = "  -9.78%  "

$ws.Range("D43").Value = "5.986"
$ws.Range("E43").Value = "  -9.43%  "

$ws.Range("D44").Value = "7.960"
$ws.Range("E44").Value = "  -5.74%  "

$ws.Range("D45").Value = "1.011"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").Value = "104.75"
$ws.Range("E46").Value = "  -1.49%  "

$ws.Range("D47").Value = "10.17"
$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("D48").Value = "0.06306"
$ws.Range("E48").Value = "  -3.06%  "

$ws.Range("D49").Value = "0.4512"
$ws.Range("E49").Value = "  -3.69%  "

$ws.Range("E50").Value = "  -2.49%  "

$ws.Range("D51").Value = "1.714"
$ws.Range("E51").Value = "  -2.71%  "
